# Update the "cryptos" price/volume table with the latest scraped values.
# Cells that hold numeric-looking text (plain decimal prices) are written
# with a leading apostrophe so Excel keeps them as text instead of
# converting them to real numbers (matching the original inlineStr cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.855.25"
$ws.Range("E2").Value = "  -1.44%  "
$ws.Range("D3").Value = "2.625.91"
$ws.Range("E3").Value = "  +0.86%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'513.01"
$ws.Range("E5").Value = "  -0.22%  "
$ws.Range("D6").Value = "'143.63"
$ws.Range("E6").Value = "  -1.72%  "
$ws.Range("D7").Value = "'0.996"
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("D8").Value = "'0.570"
$ws.Range("E8").Value = "  +1.29%  "
$ws.Range("D9").Value = "2.649.14"
$ws.Range("E9").Value = "  +1.73%  "
$ws.Range("D10").Value = "'6.29"
$ws.Range("E10").Value = "  +1.01%  "
$ws.Range("E11").Value = "  +1.75%  "
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("E13").Value = "  -1.44%  "
$ws.Range("D14").Value = "3.090.58"
$ws.Range("E14").Value = "  +1.05%  "
$ws.Range("D15").Value = "58.835.58"
$ws.Range("E16").Value = "  +0.65%  "
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("D18").Value = "2.642.65"
$ws.Range("E18").Value = "  +1.34%  "
$ws.Range("D19").Value = "'4.53"
$ws.Range("E19").Value = "  -0.99%  "
$ws.Range("D20").Value = "'342.04"
$ws.Range("E20").Value = "  +1.31%  "
$ws.Range("E21").Value = "  +1.14%  "
$ws.Range("E22").Value = "  +1.01%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").Value = "'60.86"
$ws.Range("E24").Value = "  +0.43%  "
$ws.Range("D25").Value = "'0.419"
$ws.Range("E25").Value = "  +1.43%  "
$ws.Range("D26").Value = "2.758.97"
$ws.Range("E26").Value = "  +1.44%  "
$ws.Range("E28").Value = "  +2.34%  "
$ws.Range("D29").Value = "0.0₃0801"
$ws.Range("D30").Value = "'7.09"
$ws.Range("E30").Value = "  +2.24%  "
$ws.Range("D32").Value = "'6.39"
$ws.Range("E32").Value = "  +8.27%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "'1.57"
$ws.Range("E33").Value = "  +0.64%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "'18.86"
$ws.Range("E34").Value = "  +0.68%  "
$ws.Range("D35").Value = "'149.10"
$ws.Range("E35").Value = "  -0.54%  "
$ws.Range("D36").Value = "'1.01"
$ws.Range("E36").Value = "  +12.54%  "
$ws.Range("E37").Value = "  +3.74%  "
$ws.Range("E38").Value = "  +2.47%  "
$ws.Range("E39").Value = "  +1.29%  "
$ws.Range("D40").Value = "'36.48"
$ws.Range("E40").Value = "  -0.24%  "
$ws.Range("E41").Value = "  +3.14%  "
$ws.Range("E42").Value = "  -0.25%  "
$ws.Range("D43").Value = "'280.33"
$ws.Range("E43").Value = "  -1.40%  "
$ws.Range("D44").Value = "'0.613"
$ws.Range("E44").Value = "  -0.85%  "
$ws.Range("D45").Value = "'0.994"
$ws.Range("E45").Value = "  -0.42%  "
$ws.Range("D46").Value = "'0.0984"
$ws.Range("E46").Value = "  -0.62%  "
$ws.Range("E47").Value = "  +2.66%  "
$ws.Range("E48").Value = "  -1.07%  "
$ws.Range("D49").Value = "'10.27"
$ws.Range("E49").Value = "  -0.97%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "'4.69"
$ws.Range("E50").Value = "  +2.94%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "'0.0228"
$ws.Range("E51").Value = "  -0.74%  "
